$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.850.04'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.77%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.648.59'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('E4').Value = '  +0.56%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.68%  '
$ws.Range('E6').Value = '  -0.58%  '
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0628'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.25'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.29%  '
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.640.99'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('E13').Value = '  -0.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.529'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.83'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.820.91'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.51%  '
$ws.Range('D17').Value = '0.0₃0736'
$ws.Range('E17').Value = '  -1.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '215.01'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.10%  '
$ws.Range('E19').Value = '  +0.53%  '
$ws.Range('E20').Value = '  +0.81%  '
$ws.Range('E21').Value = '  +10.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.27'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.60%  '
$ws.Range('E23').Value = '  -1.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '147.22'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.20%  '
$ws.Range('E25').Value = '  +0.45%  '
$ws.Range('E26').Value = '  -0.96%  '
$ws.Range('E27').Value = '  +0.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.68'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0509'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.55%  '
$ws.Range('E30').Value = '  +0.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.36'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.01'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.296.61'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.48%  '
$ws.Range('E34').Value = '  -0.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.45'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.69%  '
$ws.Range('E36').Value = '  -2.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.537'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.824'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.64%  '
$ws.Range('E39').Value = '  +0.52%  '
$ws.Range('E40').Value = '  -0.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.25'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.39%  '
$ws.Range('E42').Value = '  -2.21%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.787.11'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '61.91'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '91.83'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.38%  '
$ws.Range('E46').Value = '  +1.21%  '
$ws.Range('E47').Value = '  -0.84%  '
$ws.Range('E48').Value = '  +0.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.65'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.92%  '
$ws.Range('E50').Value = '  +0.31%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.407'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.27%  '
